$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "26.862.03"
Set-TextCell "E2" "  +0.78%  "
Set-TextCell "D3" "1.641.80"
Set-TextCell "E3" "  +0.68%  "
Set-TextCell "E4" "  -0.50%  "
Set-TextCell "D5" "216.46"
Set-TextCell "E5" "  -0.53%  "
Set-TextCell "D6" "0.507"
Set-TextCell "E6" "  +2.00%  "
Set-TextCell "E7" "  -0.48%  "
Set-TextCell "E8" "  +1.88%  "
Set-TextCell "E9" "  +0.37%  "
Set-TextCell "D10" "19.78"
Set-TextCell "E10" "  +4.44%  "
Set-TextCell "E11" "  +0.51%  "
Set-TextCell "D12" "1.872.01"
Set-TextCell "D13" "1.649.29"
Set-TextCell "E13" "  +1.17%  "
Set-TextCell "E14" "  +0.76%  "
Set-TextCell "E15" "  +1.56%  "
Set-TextCell "D16" "66.30"
Set-TextCell "E16" "  +3.82%  "
Set-TextCell "D17" "26.875.02"
Set-TextCell "E17" "  +0.83%  "
Set-TextCell "D18" "0.0₃0727"
Set-TextCell "E18" "  +0.99%  "
Set-TextCell "D19" "219.36"
Set-TextCell "E19" "  +3.96%  "
Set-TextCell "E20" "  -0.56%  "
Set-TextCell "E21" "  +2.00%  "
Set-TextCell "D22" "6.62"
Set-TextCell "E22" "  +7.28%  "
Set-TextCell "E23" "  +3.45%  "
Set-TextCell "D24" "9.18"
Set-TextCell "E24" "  +0.41%  "
Set-TextCell "D25" "146.00"
Set-TextCell "E25" "  -0.35%  "
Set-TextCell "E26" "  -0.70%  "
Set-TextCell "E27" "  +5.84%  "
Set-TextCell "E28" "  +1.90%  "
Set-TextCell "D29" "15.80"
Set-TextCell "E29" "  +2.03%  "
Set-TextCell "D30" "0.0507"
Set-TextCell "E30" "  +1.33%  "
Set-TextCell "E31" "  -0.44%  "
Set-TextCell "E32" "  -0.45%  "
Set-TextCell "E33" "  +2.12%  "
Set-TextCell "D34" "1.55"
Set-TextCell "E34" "  +2.77%  "
Set-TextCell "E35" "  -0.06%  "
Set-TextCell "D36" "1.245.71"
Set-TextCell "E36" "  -0.81%  "
Set-TextCell "E37" "  +0.83%  "
Set-TextCell "D38" "0.539"
Set-TextCell "E38" "  +3.21%  "
Set-TextCell "E39" "  +3.93%  "
Set-TextCell "E40" "  -0.44%  "
Set-TextCell "D41" "0.806"
Set-TextCell "E41" "  +1.15%  "
Set-TextCell "E42" "  +2.55%  "
Set-TextCell "D43" "1.783.69"
Set-TextCell "D44" "2.09"
Set-TextCell "E44" "  -2.81%  "
Set-TextCell "D45" "60.73"
Set-TextCell "E45" "  +1.88%  "
Set-TextCell "D46" "91.47"
Set-TextCell "E46" "  +0.64%  "
Set-TextCell "E47" "  +0.86%  "
Set-TextCell "B48" "BabyDogeCoin"
Set-TextCell "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D48" "0.0₆0106"
Set-TextCell "E48" "  +14.90%  "
Set-TextCell "B49" "Cronos"
Set-TextCell "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D49" "0.0515"
Set-TextCell "E49" "  -0.15%  "
Set-TextCell "B50" "Algorand"
Set-TextCell "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D50" "0.0975"
Set-TextCell "E50" "  +2.43%  "
Set-TextCell "B51" "EnergySwap"
Set-TextCell "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D51" "7.60"
Set-TextCell "E51" "  +1.76%  "
